$d = $word.ActiveDocument

# Locate the paragraph that ends with "Mikael" (the bookmark-holding paragraph)
$mikaelPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Mikael*") {
        $mikaelPara = $p
        break
    }
}

# Insert a new paragraph right after the Mikael paragraph, containing the quote.
$newRange = $mikaelPara.Range.InsertParagraphAfter()

# $newRange now represents the paragraph mark that was inserted; get the new paragraph's range
$quotePara = $mikaelPara.Next()
$quoteRange = $quotePara.Range

# Set formatting on the quote run/paragraph mark before inserting text
$quoteRange.Font.Name = "Helvetica"
$quoteRange.Font.Size = 9
$quoteRange.Font.Color = 4939595   # BGR value of 0x4B4F56 -> matches RGB() helper below
$quoteRange.Shading.BackgroundPatternColor = 15790321  # fill F1F0F0

$quoteRange.InsertBefore('"clefs : on les trouve par terre, sur des pieds d''estal" Lucas 2017')
